# Users.xlsx v1.2 update
# - Insert a new "EndTime" column between the existing "password" column (K)
#   and the trailing "don't remove" helper column (old L, which shifts to M).
# - Fill in sample data: userId "alit", password "#Parrsoo2020#",
#   EndTime "1400/10/20 13:13:13.259".
# - Resize the affected columns and move the active selection to D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column at position 12 (L); existing L (and everything to its
# right) shifts one column over to M, carrying its formatting with it.
$ws.Columns.Item(12).Insert()

# New column header (row 1) -- first newly-introduced shared string.
$ws.Range("L1").Value = "EndTime"

# Sample row (row 2) updates, entered in the same order the strings were
# first introduced so shared-string indices line up:
#   K2 -> password sample
#   L2 -> EndTime sample
#   A2 -> userId sample
$ws.Range("K2").Value = "#Parrsoo2020#"
$ws.Range("L2").Value = "1400/10/20 13:13:13.259"
$ws.Range("A2").Value = "alit"

# Column width tweaks for the new/shifted columns K, L, M.
$ws.Columns.Item(11).ColumnWidth = 11.65
$ws.Columns.Item(12).ColumnWidth = 20.79
$ws.Columns.Item(13).ColumnWidth = 44.79

# Move the selected/active cell.
$ws.Range("D10").Select() | Out-Null
